$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 189
$wsOff.Range("C3").Value = 134
$wsOff.Range("D3").Value = 48
$wsOff.Range("E3").Value = 27
$wsOff.Range("F3").Value = 5

# Sheet "DEF" - row 3 updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 216
$wsDef.Range("C3").Value = 153
$wsDef.Range("D3").Value = 71
$wsDef.Range("E3").Value = 39
$wsDef.Range("F3").Value = 4
$wsDef.Range("G3").Value = 3
